$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.531.03'
$ws.Range('E2').Value = '  +1.28%  '
$ws.Range('D3').Value = '2.987.91'
$ws.Range('E3').Value = '  +2.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '381.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.62%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.50'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.17%  '
$ws.Range('E7').Value = '  +1.58%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  +2.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.31'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.66%  '
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0848'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.07%  '
$ws.Range('D13').Value = '3.455.67'
$ws.Range('E13').Value = '  +2.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.43'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.58'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.01%  '
$ws.Range('D16').Value = '2.988.63'
$ws.Range('E16').Value = '  +2.88%  '
$ws.Range('E17').Value = '  +5.79%  '
$ws.Range('D18').Value = '51.482.65'
$ws.Range('E18').Value = '  +1.21%  '
$ws.Range('E19').Value = '  +4.07%  '
$ws.Range('E20').Value = '  +4.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.05%  '
$ws.Range('E22').Value = '  +2.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '263.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.44%  '
$ws.Range('E26').Value = '  +16.88%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.73'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +24.14%  '
$ws.Range('E28').Value = '  +14.64%  '
$ws.Range('E29').Value = '  +3.07%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.91'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '34.87'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.55%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '50.99'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('B35').Value = 'Toncoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.07'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0451'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.81%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('E38').Value = '  +2.88%  '
$ws.Range('E39').Value = '  +1.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.60'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.34%  '
$ws.Range('E41').Value = '  +0.91%  '
$ws.Range('E42').Value = '  +4.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '122.49'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.86%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.86'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.279'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +19.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.80%  '
$ws.Range('E47').Value = '  +2.87%  '
$ws.Range('E48').Value = '  +5.01%  '
$ws.Range('D49').Value = '2.031.47'
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0332'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +8.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '58.23'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.83%  '
